# Add the new "Iqaluit" weather-file row (row 81) to the WeatherData1 sheet.
# New row 82 (the trailing blank spacer row that used to be row 82) shifts
# down automatically because we are writing into row 81 which previously
# held that blank spacer.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WeatherData1")

# Values are written in the same order the original author entered them
# (file, energy_plus_location_name, country, state/province, city, ...,
# swh_fueltype, then finally location_name) so that new shared-string
# entries land in the same order as the source workbook.
$ws.Range("A81").Value = "CAN_NU_Iqaluit.AP.719090_CWEC2016.epw"
$ws.Range("C81").Value = "Iqaluit AP_NU_CAN"
$ws.Range("D81").Value = "  'CAN'"
$ws.Range("E81").Value = "NU"
$ws.Range("F81").Value = "Iqaluit"
$ws.Range("G81").Value = 9794
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 63.75
$ws.Range("K81").Value = -68.55
$ws.Range("L81").Value = 34
$ws.Range("M81").Value = 11
$ws.Range("N81").Value = 8
$ws.Range("O81").Value = "FuelOilNo2"
$ws.Range("P81").Value = "  'Hot Water'"
$ws.Range("Q81").Value = "  true"
$ws.Range("R81").Value = "  'Hot Water'"
$ws.Range("S81").Value = "  'DX'"
$ws.Range("T81").Value = "  'Scroll'"
$ws.Range("U81").Value = "  'Electric'"
$ws.Range("V81").Value = "  'Electric'"
$ws.Range("W81").Value = "  'Electric'"
$ws.Range("X81").Value = "  'var_speed_drive' "
$ws.Range("Y81").Value = "  'Electricity'"
$ws.Range("B81").Value = "CAN_NU_Iqaluit"

# Move the active cell to B84, matching the saved window/selection state.
$ws.Range("B84").Select()
